$d = $word.ActiveDocument

# Locate the paragraph that ends the "Reactive / Event Driven..." bullet -
# the three new bullets get inserted right after it, before the existing
# "Pattern Matching..." bullets.
$anchorRange = $d.Content
$found = $anchorRange.Find.Execute(
    "Reactive / Event Driven: Verticles DIDs (Distributed IDs) distributed patterns routing registry. Resource / Applicable graph logs. Rx Facade. Resource URNs Verticle Resolution, Transforms ordered Mappings Statements.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph 'Reactive / Event Driven...' not found"
}

$anchorParagraph = $anchorRange.Paragraphs(1)

$newBullets = @(
    "Statements Cases: one for each Statement - CSPO destructuring cases. One for each CSPOs Data Aggregation.",
    "Kinds Cases: One for each Kind Type Data / Mappings Agreggation.",
    "Transforms Cases: One for each Mapping Instance Data Aggregation Function."
)

$currentParagraph = $anchorParagraph
foreach ($bulletText in $newBullets) {
    # InsertParagraphAfter clones the current paragraph's (and its run's)
    # formatting onto the freshly created, empty paragraph.
    $currentParagraph.Range.InsertParagraphAfter()
    $nextIndex = $currentParagraph.Index + 1
    $currentParagraph = $d.Paragraphs($nextIndex)
    $currentParagraph.Range.Text = $bulletText
}
